# Replace "stone" with "rock" throughout the document (both occurrences
# inside the ["stone", "scissor", "paper"] literals in the Theme 5 task).
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("stone", $true, $false, $false, $false, $false, `
              $true, 1, $false, "rock", 2)
